$d = $word.ActiveDocument
$enDash = [char]0x2013

# -----------------------------------------------------------------------
# 1) "Now the bot id is - playitabhi@gmail.com" (mailto hyperlink):
#    drop the hyperlink + its "Internet Link" character style (on both
#    the run and the paragraph mark) and swap in the new address as
#    plain text.
# -----------------------------------------------------------------------
$h = $d.Hyperlinks(1)
$h.Delete()

$p2 = $d.Paragraphs(2)
$fullP2 = $d.Range($p2.Range.Start, $p2.Range.End)
$fullP2.Delete()

$newP2 = $d.Paragraphs(2)
$insertPoint2 = $d.Range($newP2.Range.Start, $newP2.Range.Start)
$insertPoint2.InsertParagraphBefore()

$targetP2 = $d.Paragraphs(2)
$r2 = $d.Range($targetP2.Range.Start, $targetP2.Range.Start)
$r2.InsertAfter("Now the bot id is " + $enDash + " demolinkedin@webaroo.com")

# -----------------------------------------------------------------------
# 2) "1. link -  a url is given ..." -> "1. help -  an image/url is given ..."
# -----------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$full4 = $p4.Range.Text
$search4 = "1. link -  a url is given"
$idx4 = $full4.IndexOf($search4)
$subStart4 = $p4.Range.Start + $idx4
$subEnd4 = $subStart4 + $search4.Length
$subRng4 = $d.Range($subStart4, $subEnd4)
$subRng4.Text = "1. help -  an image/url is given"

# -----------------------------------------------------------------------
# 3) "4. share - post on linked in" -> "4. post - post on linked in"
# -----------------------------------------------------------------------
$p6 = $d.Paragraphs(6)
$full6 = $p6.Range.Text
$search6 = "share"
$idx6 = $full6.IndexOf($search6)
$subStart6 = $p6.Range.Start + $idx6
$subEnd6 = $subStart6 + $search6.Length
$subRng6 = $d.Range($subStart6, $subEnd6)
$subRng6.Text = "post"

# -----------------------------------------------------------------------
# 4) New paragraph "coming soon -" right after the "... post on linked in"
#    line, before "5. company ...".
# -----------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$insertPoint7 = $d.Range($p7.Range.Start, $p7.Range.Start)
$insertPoint7.InsertParagraphBefore()

$newP7 = $d.Paragraphs(7)
$r7 = $d.Range($newP7.Range.Start, $newP7.Range.Start)
$r7.InsertAfter("coming soon -")

Write-Output "done"
